$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
try {
  $pc = $wb.PivotCaches().Create(1, "Journal")
  Write-Host "pc created"
  $pt = $pc.CreatePivotTable("H116", "TestPivot")
  Write-Host ("Created pivot table: " + $pt.Name)
} catch {
  Write-Host ("ERROR1: " + $_.Exception.Message)
}
